{"js": "// \"start renaming classification to label\"\n//\n// Three changes, mirroring the author's (in-progress) rename of\n// \"classification\"/\"class\" to \"label\":\n//   1. In the Document properties list: \"class\" -> \"label\"\n//   2. In the Page properties list: \"class\" -> \"label\", and the word\n//      \"classification\" gets split into \"classific\" | \"ation\" with the\n//      document's \"_GoBack\" bookmark landing around \"classific\" (this is\n//      simply where Word's internal last-edit-position bookmark ends up\n//      after the in-place retyping).\n//   3. The \"_GoBack\" bookmark that used to sit after \"Pre-Gap (Epipog) \"\n//      is gone (it moved to the spot described above - a document can\n//      only have one \"_GoBack\" bookmark).\n\nconst body = context.document.body;\n\n// Step 0: the \"_GoBack\" bookmark must be unique in the document, so drop\n// its old location before we re-create it near \"classific\" below.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// Step 1: locate the two target paragraphs by their exact current text.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet documentPropertyPara = null;\nlet pagePropertyPara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t === \"class \u2013 The document classification.\") {\n    documentPropertyPara = paragraphs.items[i];\n  } else if (t === \"class \u2013 The page classification \") {\n    pagePropertyPara = paragraphs.items[i];\n  }\n}\nif (!documentPropertyPara || !pagePropertyPara) {\n  throw new Error(\n    \"Could not locate target paragraphs (documentPropertyPara=\" +\n      !!documentPropertyPara +\n      \", pagePropertyPara=\" +\n      !!pagePropertyPara +\n      \")\"\n  );\n}\n\n// Hunk 1: \"class\" -> \"label\" in the Document properties bullet.\nconst documentClassWord = documentPropertyPara.search(\"class\", {\n  matchCase: true,\n  matchWholeWord: true\n});\ndocumentClassWord.load(\"items\");\nawait context.sync();\ndocumentClassWord.items[0].insertText(\"label\", \"Replace\");\nawait context.sync();\n\n// Hunk 2: \"class\" -> \"label\" in the Page properties bullet.\nconst pageClassWord = pagePropertyPara.search(\"class\", {\n  matchCase: true,\n  matchWholeWord: true\n});\npageClassWord.load(\"items\");\nawait context.sync();\npageClassWord.items[0].insertText(\"label\", \"Replace\");\nawait context.sync();\n\n// Hunk 2 (continued): give \" - The page \" its own run, separate from\n// \"label\" before it and \"classific\"/\"ation\" after it.\nconst pageLead = pagePropertyPara.search(\" \u2013 The page \", { matchCase: true });\npageLead.load(\"items\");\nawait context.sync();\npageLead.items[0].insertBookmark(\"__tmp_split_marker__\");\nawait context.sync();\ncontext.document.deleteBookmark(\"__tmp_split_marker__\");\nawait context.sync();\n\n// Hunk 2 (continued): re-create \"_GoBack\" around \"classific\", splitting\n// \"classification\" into \"classific\" + \"ation\".\nconst classificStem = pagePropertyPara.search(\"classific\", {\n  matchCase: true\n});\nclassificStem.load(\"items\");\nawait context.sync();\nclassificStem.items[0].insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# \"start renaming classification to label\"\n#\n# Three changes, mirroring the author's (in-progress) rename of\n# \"classification\"/\"class\" to \"label\":\n#   1. In the Document properties list: \"class\" -> \"label\"\n#   2. In the Page properties list: \"class\" -> \"label\", and the word\n#      \"classification\" gets split into \"classific\" | \"ation\" with the\n#      document's \"_GoBack\" bookmark landing around \"classific\" (this is\n#      simply where Word's internal last-edit-position bookmark ends up\n#      after the in-place retyping).\n#   3. The \"_GoBack\" bookmark that used to sit after \"Pre-Gap (Epipog) \"\n#      is gone (it moved to the spot described above - a document can\n#      only have one \"_GoBack\" bookmark).\n\n$d = $word.ActiveDocument\n$dash = [char]8211\n\n# Step 0: \"_GoBack\" must stay unique, so drop its old location before\n# re-creating it near \"classific\" further down.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n  $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# Step 1: locate the two target paragraphs by their exact current text.\n$targetText1 = \"class \" + $dash + \" The document classification.`r\"\n$targetText2 = \"class \" + $dash + \" The page classification `r\"\n\n$documentPropertyPara = $null\n$pagePropertyPara = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n  $p = $d.Paragraphs.Item($i)\n  $t = $p.Range.Text\n  if ($t -eq $targetText1) {\n    $documentPropertyPara = $p\n  } elseif ($t -eq $targetText2) {\n    $pagePropertyPara = $p\n  }\n}\nif ($null -eq $documentPropertyPara) { throw \"Could not locate the Document properties 'class' paragraph\" }\nif ($null -eq $pagePropertyPara) { throw \"Could not locate the Page properties 'class' paragraph\" }\n\n# Hunk 1: \"class\" -> \"label\" in the Document properties bullet. Wrap the\n# replace in a throwaway bookmark split so the runs that follow \"class\"\n# (\" - \", \"T\", \"he document classification\", \".\") are left untouched.\n$rng1 = $documentPropertyPara.Range\n$find1 = $rng1.Find\n$find1.ClearFormatting()\n$find1.Text = \"class\"\n$find1.MatchWholeWord = $true\n$find1.MatchCase = $true\n$find1.Forward = $true\n$find1.Wrap = 0\nif (-not $find1.Execute()) { throw \"Could not find 'class' in the Document properties paragraph\" }\n\n$splitPoint1 = $d.Range($rng1.End, $rng1.End)\n$d.Bookmarks.Add(\"__tmp_split_marker_1__\", $splitPoint1)\n$rng1.Text = \"label\"\n$d.Bookmarks.Item(\"__tmp_split_marker_1__\").Delete()\n\n# Hunk 2: Page properties bullet.\n# Step A: wrap \"classific\" with the relocated \"_GoBack\" bookmark first -\n# doing this before the \"class\" -> \"label\" replace keeps the surrounding\n# run whitespace handling clean (matches how Word itself would leave it).\n$rngA = $pagePropertyPara.Range\n$findA = $rngA.Find\n$findA.ClearFormatting()\n$findA.Text = \"classific\"\n$findA.MatchCase = $true\n$findA.Forward = $true\n$findA.Wrap = 0\nif (-not $findA.Execute()) { throw \"Could not find 'classific' in the Page properties paragraph\" }\n$d.Bookmarks.Add(\"_GoBack\", $rngA)\n\n# Step B: \"class\" -> \"label\" in the same bullet.\n$rngB = $pagePropertyPara.Range\n$findB = $rngB.Find\n$findB.ClearFormatting()\n$findB.Text = \"class\"\n$findB.MatchWholeWord = $true\n$findB.MatchCase = $true\n$findB.Forward = $true\n$findB.Wrap = 0\nif (-not $findB.Execute()) { throw \"Could not find 'class' in the Page properties paragraph\" }\n\n$splitPoint2 = $d.Range($rngB.End, $rngB.End)\n$d.Bookmarks.Add(\"__tmp_split_marker_2__\", $splitPoint2)\n$rngB.Text = \"label\"\n$d.Bookmarks.Item(\"__tmp_split_marker_2__\").Delete()\n"}
